# Update the two-digit division table cells to the new values.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Text = "89÷8="
$cell = $t.Cell(1, 2)
$cell.Range.Text = "11÷9="
$cell = $t.Cell(1, 3)
$cell.Range.Text = "23÷8="
$cell = $t.Cell(1, 4)
$cell.Range.Text = "91÷8="
$cell = $t.Cell(1, 5)
$cell.Range.Text = "19÷8="

$cell = $t.Cell(5, 1)
$cell.Range.Text = "89÷6="
$cell = $t.Cell(5, 2)
$cell.Range.Text = "89÷3="
$cell = $t.Cell(5, 3)
$cell.Range.Text = "56÷7="
$cell = $t.Cell(5, 4)
$cell.Range.Text = "19÷5="
$cell = $t.Cell(5, 5)
$cell.Range.Text = "81÷5="

$cell = $t.Cell(9, 1)
$cell.Range.Text = "31÷8="
$cell = $t.Cell(9, 2)
$cell.Range.Text = "25÷3="
$cell = $t.Cell(9, 3)
$cell.Range.Text = "50÷9="
$cell = $t.Cell(9, 4)
$cell.Range.Text = "56÷9="
$cell = $t.Cell(9, 5)
$cell.Range.Text = "52÷3="

$cell = $t.Cell(13, 1)
$cell.Range.Text = "49÷9="
$cell = $t.Cell(13, 2)
$cell.Range.Text = "57÷2="
$cell = $t.Cell(13, 3)
$cell.Range.Text = "65÷2="
$cell = $t.Cell(13, 4)
$cell.Range.Text = "31÷3="
$cell = $t.Cell(13, 5)
$cell.Range.Text = "41÷4="

$cell = $t.Cell(17, 1)
$cell.Range.Text = "54÷2="
$cell = $t.Cell(17, 2)
$cell.Range.Text = "64÷6="
$cell = $t.Cell(17, 3)
$cell.Range.Text = "85÷7="
$cell = $t.Cell(17, 4)
$cell.Range.Text = "76÷5="
$cell = $t.Cell(17, 5)
$cell.Range.Text = "41÷9="
